$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 250; this shifts the existing rows 250..318
# down to 251..319 and extends the used range to A1:R319.
$ws.Range("A250").EntireRow.Insert()

# Populate the newly inserted row 250 with the new record.
$ws.Range("A250").Value = 4
$ws.Range("B250").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C250").Value = "Los Lagos"
$ws.Range("D250").Value = 44736
$ws.Range("E250").Value = 10
$ws.Range("F250").Value = 100112045
$ws.Range("G250").Value = "Zapallo"
$ws.Range("H250").Value = "Paine"
$ws.Range("I250").Value = "1a (guarda)"
$ws.Range("J250").Value = 1000
$ws.Range("K250").Value = 500
$ws.Range("L250").Value = 500
$ws.Range("M250").Value = 500
$ws.Range("N250").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O250").Value = "Región de O'Higgins"
$ws.Range("P250").Value = 500
$ws.Range("Q250").Value = 1
$ws.Range("R250").Value = "Hortaliza"
